$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as plain text (not numbers) so that values like
# "1.80" or "0.999" keep their exact original formatting/decimal places
# instead of Excel auto-converting numeric-looking text into a real number.
# Force text format on the specific D cells whose new value would otherwise
# be auto-detected as a number.

# Row 2
$ws.Range("D2").Value = "47.287.23"
$ws.Range("E2").Value = "  +5.98%  "

# Row 3
$ws.Range("D3").Value = "2.508.36"
$ws.Range("E3").Value = "  +3.72%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.66"
$ws.Range("E5").Value = "  +2.81%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.05"
$ws.Range("E6").Value = "  +4.80%  "

# Row 7
$ws.Range("E7").Value = "  +2.20%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.08%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.541"
$ws.Range("E9").Value = "  +2.50%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.64"
$ws.Range("E10").Value = "  +3.90%  "

# Row 11
$ws.Range("E11").Value = "  +2.42%  "

# Row 12
$ws.Range("E12").Value = "  +0.88%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.38"
$ws.Range("E13").Value = "  -1.17%  "

# Row 14
$ws.Range("E14").Value = "  +4.02%  "

# Row 15
$ws.Range("D15").Value = "2.897.93"
$ws.Range("E15").Value = "  +3.62%  "

# Row 16
$ws.Range("D16").Value = "2.549.89"
$ws.Range("E16").Value = "  +6.53%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.847"
$ws.Range("E17").Value = "  +2.30%  "

# Row 18
$ws.Range("D18").Value = "47.161.97"
$ws.Range("E18").Value = "  +6.01%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.76"
$ws.Range("E19").Value = "  +4.72%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.52"
$ws.Range("E20").Value = "  +2.79%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0942"
$ws.Range("E21").Value = "  +2.93%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.95"
$ws.Range("E22").Value = "  +3.45%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.98"
$ws.Range("E23").Value = "  +4.65%  "

# Row 24
$ws.Range("E24").Value = "  +5.57%  "

# Row 25
$ws.Range("E25").Value = "  +2.82%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.44"
$ws.Range("E26").Value = "  +5.00%  "

# Row 27
$ws.Range("E27").Value = "  +0.01%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.95"
$ws.Range("E28").Value = "  +4.83%  "

# Row 29
$ws.Range("E29").Value = "  -3.40%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.28"
$ws.Range("E30").Value = "  +6.06%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.137"
$ws.Range("E31").Value = "  +8.94%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.78"
$ws.Range("E32").Value = "  +3.08%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.81"
$ws.Range("E33").Value = "  +1.61%  "

# Row 34
$ws.Range("E34").Value = "  +3.39%  "

# Row 35
$ws.Range("E35").Value = "  +1.56%  "

# Row 36
$ws.Range("E36").Value = "  -0.03%  "

# Row 37
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.65"
$ws.Range("E37").Value = "  +4.48%  "

# Row 38
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.94"
$ws.Range("E38").Value = "  +3.44%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.96"
$ws.Range("E39").Value = "  +4.29%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "124.48"
$ws.Range("E40").Value = "  -0.19%  "

# Row 41
$ws.Range("E41").Value = "  +2.49%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.23"
$ws.Range("E42").Value = "  +2.02%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.25"
$ws.Range("E43").Value = "  +2.15%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0297"
$ws.Range("E44").Value = "  +2.87%  "

# Row 45
$ws.Range("D45").Value = "1.981.40"
$ws.Range("E45").Value = "  +2.13%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.03"
$ws.Range("E46").Value = "  +3.51%  "

# Row 47
$ws.Range("E47").Value = "  -0.04%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.80"
$ws.Range("E48").Value = "  +4.07%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.41"
$ws.Range("E49").Value = "  +17.91%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.05"
$ws.Range("E50").Value = "  -1.06%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.07"
$ws.Range("E51").Value = "  +6.66%  "
